# Updated cryptos list on Wed Jun 26 17:24:06 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a price-column value while keeping it stored as plain text
# (matches source data, which always stores Price cells as inline strings,
# never numbers) instead of letting Excel auto-coerce numeric-looking text
# into a real number. Style is reset to Normal afterwards so no stray
# cell-format attribute is left behind on cells that never had one.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "60.907.29"
$ws.Range("E2").Value = "  -1.18%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.338.49"
$ws.Range("E3").Value = "  -1.58%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "568.44"
$ws.Range("E5").Value = "  -1.19%  "

# Row 6 - Solana
Set-TextValue "D6" "135.03"
$ws.Range("E6").Value = "  -2.18%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - LidoStakedEther
Set-TextValue "D8" "3.335.55"
$ws.Range("E8").Value = "  -1.63%  "

# Row 9 - XRP
Set-TextValue "D9" "0.467"
$ws.Range("E9").Value = "  -2.46%  "

# Row 10 - Toncoin
Set-TextValue "D10" "7.42"
$ws.Range("E10").Value = "  -1.52%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.122"
$ws.Range("E11").Value = "  -3.56%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.385"
$ws.Range("E12").Value = "  -2.68%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "3.907.87"
$ws.Range("E13").Value = "  -1.51%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.46%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  -4.32%  "

# Row 16 - Avalanche
Set-TextValue "D16" "25.70"
$ws.Range("E16").Value = "  +1.37%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "3.336.98"
$ws.Range("E17").Value = "  -2.03%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "61.055.10"
$ws.Range("E18").Value = "  -1.05%  "

# Row 19 - Chainlink
Set-TextValue "D19" "13.93"

# Row 20 - Polkadot
Set-TextValue "D20" "5.80"
$ws.Range("E20").Value = "  -1.76%  "

# Row 21 - Uniswap
Set-TextValue "D21" "9.19"
$ws.Range("E21").Value = "  -1.90%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "376.98"
$ws.Range("E22").Value = "  -2.94%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.549"
$ws.Range("E23").Value = "  -4.12%  "

# Row 24 - WrappedeETH
Set-TextValue "D24" "3.491.08"
$ws.Range("E24").Value = "  -1.03%  "

# Row 25 - Dai
Set-TextValue "D25" "1.00"
$ws.Range("E25").Value = "  +0.06%  "

# Row 26 - Litecoin
Set-TextValue "D26" "70.67"
$ws.Range("E26").Value = "  -0.55%  "

# Row 27 - PEPE
Set-TextValue "D27" "0.0000122"
$ws.Range("E27").Value = "  -4.69%  "

# Row 28 - Fetch.AI
Set-TextValue "D28" "1.74"
$ws.Range("E28").Value = "  +7.68%  "

# Row 29 & 30 swapped: RenderToken <-> Binance-PegBSC-USD (with updated data)
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  +0.49%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D30" "7.40"
$ws.Range("E30").Value = "  -4.17%  "

# Row 31 - Kaspa
Set-TextValue "D31" "0.163"
$ws.Range("E31").Value = "  +3.16%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextValue "D32" "8.11"
$ws.Range("E32").Value = "  -2.67%  "

# Row 33 - PancakeSwap
Set-TextValue "D33" "2.12"
$ws.Range("E33").Value = "  -1.48%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  +0.04%  "

# Row 35 - EthereumClassic
Set-TextValue "D35" "23.38"
$ws.Range("E35").Value = "  -0.55%  "

# Row 36 - NEARProtocol
Set-TextValue "D36" "5.18"
$ws.Range("E36").Value = "  -6.27%  "

# Row 37 - Aptos
Set-TextValue "D37" "6.73"
$ws.Range("E37").Value = "  -3.77%  "

# Row 38 - Monero
Set-TextValue "D38" "164.62"
$ws.Range("E38").Value = "  +1.57%  "

# Row 39 - ImmutableX
Set-TextValue "D39" "1.51"
$ws.Range("E39").Value = "  -2.59%  "

# Row 40 - Hedera
Set-TextValue "D40" "0.0754"
$ws.Range("E40").Value = "  -5.21%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  -0.05%  "

# Row 42 - Mantle
Set-TextValue "D42" "0.764"
$ws.Range("E42").Value = "  -1.09%  "

# Row 43 & 44 swapped: Stacks <-> OKB (with updated data)
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D43" "41.30"
$ws.Range("E43").Value = "  -0.16%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D44" "1.69"
$ws.Range("E44").Value = "  -2.91%  "

# Row 45 - ONDO
Set-TextValue "D45" "1.19"
$ws.Range("E45").Value = "  -2.84%  "

# Row 46 - Filecoin
Set-TextValue "D46" "4.35"
$ws.Range("E46").Value = "  -2.44%  "

# Row 47 - EnergySwap
Set-TextValue "D47" "23.64"
$ws.Range("E47").Value = "  -4.40%  "

# Row 48 - Cosmos
Set-TextValue "D48" "6.76"
$ws.Range("E48").Value = "  -3.19%  "

# Row 49 - InjectiveProtocol
Set-TextValue "D49" "22.71"
$ws.Range("E49").Value = "  -1.50%  "

# Row 50 - Maker
Set-TextValue "D50" "2.338.04"
$ws.Range("E50").Value = "  -1.91%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  -2.52%  "
